$wb = $excel.ActiveWorkbook

# --- Rename the original sheet to "Sheet1" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sheet1"

# --- Update a few cells on Sheet1 (Import page tweaks) ---
$ws1.Range("B2").Value = "This"
$ws1.Range("B4").Value = "Sheet"
$ws1.Range("B5").Value = "Number "
$ws1.Range("B6").Value = "One"

# --- Add a second sheet right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Populate Sheet2 with the same base table, customized for Sheet2 ---
$ws2.Range("A1").Value = "Col1"
$ws2.Range("B1").Value = "Col2"
$ws2.Range("C1").Value = "Col3"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "This "
$ws2.Range("C2").Value = "T"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "Is "
$ws2.Range("C3").Value = "F"

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = "Sheet"
$ws2.Range("C4").Value = "T"

$ws2.Range("A5").Value = 4
$ws2.Range("B5").Value = "Number"
$ws2.Range("C5").Value = "F"

$ws2.Range("A6").Value = 5
$ws2.Range("B6").Value = "Two"
$ws2.Range("C6").Value = "F"

# --- Selections matching the final state ---
[void]$ws2.Range("F18").Select()
[void]$ws1.Activate()
[void]$ws1.Range("B6").Select()
